$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether | Ether
$ws.Range("H15").Value = 1385.5393
$ws.Range("I15").Value = 1385.5393
$ws.Range("K15").Value = 4156.617899999999
$ws.Range("M15").Value = -3987.617899999999

# Row 26: Everything Is Impossible | Budding Ash Wand
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 45457230
$ws.Range("I62").Value = 62501830
$ws.Range("J62").Value = 4968.3335
$ws.Range("K62").Value = 62501830
$ws.Range("L62").Value = 4968.3335
$ws.Range("M62").Value = -62501206
$ws.Range("N62").Value = -6216.3335

# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 45457230
$ws.Range("I65").Value = 62501830
$ws.Range("J65").Value = 4968.3335
$ws.Range("K65").Value = 312509150
$ws.Range("L65").Value = 24841.6675
$ws.Range("M65").Value = -312506030
$ws.Range("N65").Value = -31081.6675

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 5264.8276
$ws.Range("I132").Value = 1803.3334
$ws.Range("J132").Value = 21880
$ws.Range("K132").Value = 5410.0002
$ws.Range("L132").Value = 65640
$ws.Range("M132").Value = -2880.0002
$ws.Range("N132").Value = -70700

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Range("H133").Value = 48035
$ws.Range("J133").Value = 48035
$ws.Range("L133").Value = 48035
$ws.Range("N133").Value = -58155

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 795126
$ws.Range("I137").Value = 1193.7675
$ws.Range("J137").Value = 2059536.6
$ws.Range("K137").Value = 3581.3025
$ws.Range("L137").Value = 6178609.800000001
$ws.Range("M137").Value = -1031.3025
$ws.Range("N137").Value = -6183709.800000001

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2967.4614
$ws.Range("I138").Value = 2167.5122
$ws.Range("J138").Value = 4334.0415
$ws.Range("K138").Value = 6502.5366
$ws.Range("L138").Value = 13002.1245
$ws.Range("M138").Value = -1362.5366
$ws.Range("N138").Value = -23282.1245

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 13378.81
$ws.Range("I32").Value = 5683.1885
$ws.Range("J32").Value = 30507.773
$ws.Range("K32").Value = 5683.1885
$ws.Range("L32").Value = 30507.773
$ws.Range("M32").Value = -5396.1885
$ws.Range("N32").Value = -31081.773

# Row 43: They've Got Legs | Steel Sabatons
$ws.Range("H43").Value = 7150.4443
$ws.Range("J43").Value = 7150.4443
$ws.Range("L43").Value = 7150.4443
$ws.Range("N43").Value = -7776.4443

# Row 82: Belle of the Brawl | Titanium Vambraces of Fending
$ws.Range("H82").Value = 28000
$ws.Range("J82").Value = 28000
$ws.Range("L82").Value = 28000
$ws.Range("N82").Value = -28722

# Row 85: Shouldering the Shut-ins (L) | Titanium Vambraces of Fending
$ws.Range("H85").Value = 28000
$ws.Range("J85").Value = 28000
$ws.Range("L85").Value = 28000
$ws.Range("N85").Value = -30496

# Row 107: Shielding the Realm | Deepgold Kite Shield
$ws.Range("H107").Value = 28788
$ws.Range("J107").Value = 28788
$ws.Range("L107").Value = 28788
$ws.Range("N107").Value = -36468

$ws = $wb.Worksheets.Item("BSM")
# Row 92: Have Blade, Will Travel | High Steel Katzbalger
$ws.Range("H92").Value = 23540
$ws.Range("J92").Value = 23540
$ws.Range("L92").Value = 23540
$ws.Range("N92").Value = -28532

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2686.673
$ws.Range("I31").Value = 1139.0454
$ws.Range("J31").Value = 3821.6
$ws.Range("K31").Value = 1139.0454
$ws.Range("L31").Value = 3821.6
$ws.Range("M31").Value = -844.0454
$ws.Range("N31").Value = -4411.6

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2686.673
$ws.Range("I34").Value = 1139.0454
$ws.Range("J34").Value = 3821.6
$ws.Range("K34").Value = 1139.0454
$ws.Range("L34").Value = 3821.6
$ws.Range("M34").Value = -937.0454
$ws.Range("N34").Value = -4225.6

# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 85751.586
$ws.Range("J141").Value = 91729
$ws.Range("L141").Value = 91729
$ws.Range("N141").Value = -102089

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry | Kukuru Powder
$ws.Range("H14").Value = 7.090909
$ws.Range("I14").Value = 7.090909
$ws.Range("K14").Value = 21.272727
$ws.Range("M14").Value = 151.727273

# Row 63: The Next to Last Supper | Stuffed Cabbage Rolls
$ws.Range("H63").Value = 4933.769
$ws.Range("J63").Value = 6683
$ws.Range("L63").Value = 20049
$ws.Range("N63").Value = -21547

# Row 66: Nostalgia through the Stomach (L) | Stuffed Cabbage Rolls
$ws.Range("H66").Value = 4933.769
$ws.Range("J66").Value = 6683
$ws.Range("L66").Value = 60147
$ws.Range("N66").Value = -67635

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 974.04
$ws.Range("I131").Value = 210
$ws.Range("J131").Value = 989.6326
$ws.Range("K131").Value = 630
$ws.Range("L131").Value = 2968.8978
$ws.Range("M131").Value = 4410
$ws.Range("N131").Value = -13048.8978

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 864038.6
$ws.Range("I132").Value = 2012.4762
$ws.Range("J132").Value = 3126857.2
$ws.Range("K132").Value = 18112.2858
$ws.Range("L132").Value = 28141714.8
$ws.Range("M132").Value = -15582.2858
$ws.Range("N132").Value = -28146774.8

$ws = $wb.Worksheets.Item("GSM")
# Row 63: Not on My Table | Mythrite Earrings of Healing
$ws.Range("H63").Value = 13694.375
$ws.Range("J63").Value = 13694.375
$ws.Range("L63").Value = 13694.375
$ws.Range("N63").Value = -15066.375

# Row 66: Heinz's Dilemma (L) | Mythrite Earrings of Healing
$ws.Range("H66").Value = 13694.375
$ws.Range("J66").Value = 13694.375
$ws.Range("L66").Value = 41083.125
$ws.Range("N66").Value = -47947.125

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 5556754
$ws.Range("I102").Value = 7937489.5
$ws.Range("J102").Value = 1704.6666
$ws.Range("K102").Value = 7937489.5
$ws.Range("L102").Value = 1704.6666
$ws.Range("M102").Value = -7935867.5
$ws.Range("N102").Value = -4948.6666

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 30304174
$ws.Range("I126").Value = 30304174
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 90912522
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -90910052
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 5444.9546
$ws.Range("I132").Value = 5343.643
$ws.Range("J132").Value = 5622.25
$ws.Range("K132").Value = 16030.929
$ws.Range("L132").Value = 16866.75
$ws.Range("M132").Value = -13500.929
$ws.Range("N132").Value = -21926.75

# Row 141: Just Generally Freezing | Gargantuaskin Trousers of Striking
$ws.Range("H141").Value = 145538
$ws.Range("J141").Value = 145538
$ws.Range("L141").Value = 145538
$ws.Range("N141").Value = -155898

$ws = $wb.Worksheets.Item("WVR")
# Row 41: Half Is the New Double | Linen Halfgloves
$ws.Range("H41").Value = 7210.6
$ws.Range("J41").Value = 7210.6
$ws.Range("L41").Value = 7210.6
$ws.Range("N41").Value = -7990.6

# Row 45: Private Concerns | Linen Trousers
$ws.Range("H45").Value = 5579.273
$ws.Range("I45").Value = 3569
$ws.Range("J45").Value = 5780.3
$ws.Range("K45").Value = 3569
$ws.Range("L45").Value = 5780.3
$ws.Range("M45").Value = -3078
$ws.Range("N45").Value = -6762.3

# Row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws.Range("H74").Value = 6671
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 6671
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 6671
$ws.Range("N74").Value = -8543
$ws.Range("M74").ClearContents()

# Row 77: When in Robes (L) | Ramie Robe of Casting
$ws.Range("H77").Value = 6671
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 6671
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 20013
$ws.Range("N77").Value = -29373
$ws.Range("M77").ClearContents()
